$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 48-69: southwest plant data (first pass) ---
# row 48
$ws.Range("A48").Value = 'Prunus'
$ws.Range("B48").Value = 'serotina'
$ws.Range("C48").Value = 'syn'
$ws.Range("D48").Value = 'insect'
$ws.Range("E48").Value = 100
$ws.Range("F48").Value = 'bisexual'
$ws.Range("G48").Value = 4.5
$ws.Range("H48").Value = 8

# row 49
$ws.Range("A49").Value = 'Prunus'
$ws.Range("B49").Value = 'emarginata'
$ws.Range("C49").Value = 'syn'
$ws.Range("D49").Value = 'insect'
$ws.Range("E49").Value = 12
$ws.Range("F49").Value = 'bisexual'
$ws.Range("G49").Value = 5
$ws.Range("H49").Value = 7

# row 50
$ws.Range("A50").Value = 'Prunus'
$ws.Range("B50").Value = 'gladulosa'
$ws.Range("C50").Value = 'pro/syn'
$ws.Range("D50").Value = 'insect'
$ws.Range("E50").Value = 5
$ws.Range("F50").Value = 'bisexual'
$ws.Range("H50").Value = 4.5

# row 51
$ws.Range("A51").Value = 'Prunus'
$ws.Range("B51").Value = 'minutiflora'
$ws.Range("C51").Value = 'syn'
$ws.Range("D51").Value = 'insect'
$ws.Range("E51").Value = 3
$ws.Range("F51").Value = 'bisexual'
$ws.Range("G51").Value = 3
$ws.Range("H51").Value = 5.5

# row 52
$ws.Range("A52").Value = 'Prunus'
$ws.Range("B52").Value = 'havardii'
$ws.Range("C52").Value = 'ser'
$ws.Range("D52").Value = 'insect'
$ws.Range("F52").Value = 'bisexual'
$ws.Range("G52").Value = 5.5
$ws.Range("H52").Value = 7

# row 53
$ws.Range("A53").Value = 'Prunus'
$ws.Range("B53").Value = 'gracilis'
$ws.Range("C53").Value = 'pro'
$ws.Range("D53").Value = 'insect'
$ws.Range("E53").Value = 15
$ws.Range("F53").Value = 'bisexual'
$ws.Range("G53").Value = 3
$ws.Range("H53").Value = 7

# row 54
$ws.Range("A54").Value = 'Prunus'
$ws.Range("B54").Value = 'murrayana'
$ws.Range("C54").Value = 'syn'
$ws.Range("D54").Value = 'insect'
$ws.Range("E54").Value = 6
$ws.Range("F54").Value = 'bisexual'
$ws.Range("G54").Value = 3.5

# row 55
$ws.Range("A55").Value = 'Prunus'
$ws.Range("B55").Value = 'angustifolia'
$ws.Range("C55").Value = 'pro'
$ws.Range("D55").Value = 'insect'
$ws.Range("E55").Value = 25
$ws.Range("F55").Value = 'bisexual'
$ws.Range("G55").Value = 3.5
$ws.Range("H55").Value = 6

# row 56
$ws.Range("A56").Value = 'Prunus'
$ws.Range("B56").Value = 'reverchonii'
$ws.Range("C56").Value = 'pro/syn'
$ws.Range("D56").Value = 'insect'
$ws.Range("E56").Value = 6
$ws.Range("F56").Value = 'bisexual'
$ws.Range("G56").Value = 3.5
$ws.Range("H56").Value = 8

# row 57
$ws.Range("A57").Value = 'Prunus'
$ws.Range("B57").Value = 'rivularis'
$ws.Range("C57").Value = 'pro/syn'
$ws.Range("D57").Value = 'insect'
$ws.Range("E57").Value = 8
$ws.Range("F57").Value = 'bisexual'
$ws.Range("G57").Value = 3
$ws.Range("H57").Value = 6

# row 58
$ws.Range("A58").Value = 'Prunus'
$ws.Range("B58").Value = 'persoca'
$ws.Range("C58").Value = 'pro'
$ws.Range("D58").Value = 'insect'
$ws.Range("E58").Value = 24
$ws.Range("F58").Value = 'bisexual'
$ws.Range("G58").Value = 4
$ws.Range("H58").Value = 8.5

# row 59
$ws.Range("A59").Value = 'Pyrus'
$ws.Range("B59").Value = 'communis'
$ws.Range("C59").Value = 'pro/syn'
$ws.Range("D59").Value = 'insect'
$ws.Range("E59").Value = 75
$ws.Range("F59").Value = 'bisexual'
$ws.Range("G59").Value = 4
$ws.Range("H59").Value = 8.5

# row 60
$ws.Range("A60").Value = 'Pyrus'
$ws.Range("B60").Value = 'malus'
$ws.Range("C60").Value = 'syn'
$ws.Range("D60").Value = 'insect'
$ws.Range("E60").Value = 50
$ws.Range("F60").Value = 'bisexual'
$ws.Range("G60").Value = 4.5
$ws.Range("H60").Value = 9.5

# row 61
$ws.Range("A61").Value = 'Amelanchier'
$ws.Range("B61").Value = 'denticulata'
$ws.Range("C61").Value = 'pro/syn'
$ws.Range("D61").Value = 'insect'
$ws.Range("E61").Value = 10
$ws.Range("F61").Value = 'bisexual'
$ws.Range("G61").Value = 2.5
$ws.Range("H61").Value = 5.5

# row 62
$ws.Range("A62").Value = 'Peraphyllum'
$ws.Range("B62").Value = 'ramosissimum'
$ws.Range("C62").Value = 'syn'
$ws.Range("D62").Value = 'insect'
$ws.Range("E62").Value = 10
$ws.Range("F62").Value = 'bisexual'
$ws.Range("G62").Value = 4.5
$ws.Range("H62").Value = 8

# row 63
$ws.Range("A63").Value = 'Rubus'
$ws.Range("B63").Value = 'sons'
$ws.Range("C63").Value = 'ser'
$ws.Range("D63").Value = 'insect'
$ws.Range("E63").Value = 3
$ws.Range("F63").Value = 'bisexual'
$ws.Range("G63").Value = 5.5
$ws.Range("H63").Value = 7.5

# row 64
$ws.Range("A64").Value = 'Rubus'
$ws.Range("B64").Value = 'occidentalis'
$ws.Range("C64").Value = 'ser'
$ws.Range("D64").Value = 'insect'
$ws.Range("F64").Value = 'bisexual'
$ws.Range("G64").Value = 5.5

# row 65
$ws.Range("A65").Value = 'Rubus'
$ws.Range("B65").Value = 'nefrens'
$ws.Range("C65").Value = 'ser'
$ws.Range("D65").Value = 'insect'
$ws.Range("F65").Value = 'bisexual'
$ws.Range("G65").Value = 5.5
$ws.Range("H65").Value = 6.5

# row 66
$ws.Range("A66").Value = 'Rubus'
$ws.Range("B66").Value = 'serissimus'
$ws.Range("C66").Value = 'ser'
$ws.Range("D66").Value = 'insect'
$ws.Range("E66").Value = 4
$ws.Range("F66").Value = 'bisexual'

# row 67
$ws.Range("A67").Value = 'Rubus'
$ws.Range("B67").Value = 'allegheniensis'
$ws.Range("C67").Value = 'ser'
$ws.Range("D67").Value = 'insect'
$ws.Range("E67").Value = 3
$ws.Range("F67").Value = 'bisexual'
$ws.Range("G67").Value = 6

# row 68
$ws.Range("A68").Value = 'Prosopis'
$ws.Range("B68").Value = 'pubescens'
$ws.Range("C68").Value = 'ser'
$ws.Range("D68").Value = 'insect'
$ws.Range("E68").Value = 30
$ws.Range("F68").Value = 'bisexual'

# row 69
$ws.Range("A69").Value = 'Gymnocladus'
$ws.Range("B69").Value = 'dioica'
$ws.Range("C69").Value = 'ser'
$ws.Range("D69").Value = 'insect'
$ws.Range("E69").Value = 100
$ws.Range("F69").Value = 'bisexual'
$ws.Range("G69").Value = 5.5
$ws.Range("H69").Value = 12

# --- Row 11: add N11 (persistant) ---
$ws.Range("N11").Value = 'persistant'

# --- New rows 70-83: southwest plant data (second pass) ---
# row 70
$ws.Range("A70").Value = 'Cercis'
$ws.Range("B70").Value = 'canadensis'
$ws.Range("C70").Value = 'pro'
$ws.Range("D70").Value = 'insect'
$ws.Range("E70").Value = 40
$ws.Range("F70").Value = 'bisexual'
$ws.Range("G70").Value = 4
$ws.Range("H70").Value = 9.5

# row 71
$ws.Range("A71").Value = 'Cercis'
$ws.Range("B71").Value = 'occidentalis'
$ws.Range("C71").Value = 'pro'
$ws.Range("D71").Value = 'insect'
$ws.Range("E71").Value = 20
$ws.Range("F71").Value = 'bisexual'
$ws.Range("G71").Value = 3

# row 72
$ws.Range("A72").Value = 'Erythrina'
$ws.Range("B72").Value = 'flavelliformis'
$ws.Range("C72").Value = 'pro/syn'
$ws.Range("D72").Value = 'insect'
$ws.Range("E72").Value = 15
$ws.Range("F72").Value = 'bisexual'

# row 73
$ws.Range("A73").Value = 'Robinia'
$ws.Range("B73").Value = 'neomexicana'
$ws.Range("C73").Value = 'ser'
$ws.Range("D73").Value = 'insect'
$ws.Range("E73").Value = 24
$ws.Range("F73").Value = 'bisexual'
$ws.Range("G73").Value = 6
$ws.Range("H73").Value = 9.5

# row 74
$ws.Range("A74").Value = 'Sophora'
$ws.Range("B74").Value = 'secundiflora'
$ws.Range("C74").Value = 'syn'
$ws.Range("D74").Value = 'insect'
$ws.Range("E74").Value = 35
$ws.Range("F74").Value = 'bisexual'
$ws.Range("G74").Value = 3.5
$ws.Range("H74").Value = 9

# row 75
$ws.Range("A75").Value = 'Wisteria'
$ws.Range("B75").Value = 'frutescens'
$ws.Range("C75").Value = 'syn'
$ws.Range("D75").Value = 'insect'
$ws.Range("E75").Value = 40
$ws.Range("F75").Value = 'bisexual'
$ws.Range("G75").Value = 5
$ws.Range("H75").Value = 10

# row 76
$ws.Range("A76").Value = 'Poncirus'
$ws.Range("B76").Value = 'trifoliata'
$ws.Range("C76").Value = 'pro'
$ws.Range("D76").Value = 'insect'
$ws.Range("E76").Value = 30
$ws.Range("F76").Value = 'bisexual'
$ws.Range("G76").Value = 4.5
$ws.Range("H76").Value = 9.5

# row 77
$ws.Range("A77").Value = 'Zanthoxylem'
$ws.Range("B77").Value = 'americanum'
$ws.Range("C77").Value = 'pro'
$ws.Range("D77").Value = 'insect'
$ws.Range("E77").Value = 25
$ws.Range("F77").Value = 'dioecious'
$ws.Range("G77").Value = 4.5
$ws.Range("H77").Value = 8

# row 78
$ws.Range("A78").Value = 'Zanthoxylem'
$ws.Range("B78").Value = 'fagara'
$ws.Range("C78").Value = 'ser'
$ws.Range("D78").Value = 'insect'
$ws.Range("E78").Value = 30
$ws.Range("F78").Value = 'dioecious'
$ws.Range("G78").Value = 4.5
$ws.Range("H78").Value = 9

# row 79
$ws.Range("A79").Value = 'Zanthoxylem'
$ws.Range("B79").Value = 'parvum'
$ws.Range("C79").Value = 'pro/syn'
$ws.Range("D79").Value = 'insect'
$ws.Range("E79").Value = 4.5
$ws.Range("F79").Value = 'unisexual'

# row 80
$ws.Range("A80").Value = 'Pistacia'
$ws.Range("B80").Value = 'texana'
$ws.Range("C80").Value = 'pro/syn'
$ws.Range("D80").Value = 'insect'
$ws.Range("E80").Value = 30
$ws.Range("F80").Value = 'unisexual'

# row 81
$ws.Range("A81").Value = 'Rhus'
$ws.Range("B81").Value = 'aromatica'
$ws.Range("C81").Value = 'pro'
$ws.Range("D81").Value = 'insect'
$ws.Range("E81").Value = 8
$ws.Range("F81").Value = 'unisexual'
$ws.Range("H81").Value = 7.5

# row 82
$ws.Range("A82").Value = 'Rhus'
$ws.Range("B82").Value = 'microphylla'
$ws.Range("C82").Value = 'pro'
$ws.Range("D82").Value = 'insect'
$ws.Range("E82").Value = 15
$ws.Range("H82").Value = 6

# row 83
$ws.Range("A83").Value = 'Ilex'
$ws.Range("B83").Value = 'decidua'
$ws.Range("C83").Value = 'syn'
$ws.Range("D83").Value = 'insect'
$ws.Range("E83").Value = 30
$ws.Range("F83").Value = 'unisexual'
$ws.Range("G83").Value = 4
$ws.Range("H83").Value = 9

# --- Row 11: add O11, P11 (early autumn=9) ---
$ws.Range("O11").Value = 12
$ws.Range("P11").Value = 'early autumn=9'

# --- New row 84 ---
# row 84
$ws.Range("A84").Value = 'Ilex'
$ws.Range("B84").Value = 'montana'
$ws.Range("C84").Value = 'syn'
$ws.Range("D84").Value = 'insect'
$ws.Range("E84").Value = 40
$ws.Range("F84").Value = 'unisexual'
$ws.Range("G84").Value = 6
$ws.Range("H84").Value = 10

# --- Row 35: add D35 (insect; reuses existing shared string) ---
$ws.Range("D35").Value = 'insect'

# --- Final selection state ---
$ws.Range("A85").Select()
